$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "27.232.44"
    "D3" = "1.901.87"
    "E3" = "  +0.02%  "
    "E4" = "  +0.05%  "
    "D5" = "306.47"
    "E5" = "  -0.66%  "
    "E6" = "  -0.01%  "
    "D7" = "0.5336"
    "E7" = "  +2.28%  "
    "E8" = "  +0.87%  "
    "D9" = "0.07289"
    "E9" = "  -0.17%  "
    "D10" = "22.16"
    "E10" = "  +4.24%  "
    "D11" = "0.9024"
    "E11" = "  -0.10%  "
    "D12" = "0.08218"
    "E12" = "  +0.54%  "
    "D13" = "96.09"
    "E13" = "  -0.80%  "
    "D14" = "5.333"
    "E14" = "  +0.60%  "
    "D15" = "1.003"
    "E15" = "  +0.09%  "
    "E16" = "  +1.71%  "
    "D17" = "0.000008650"
    "E17" = "  +0.09%  "
    "E18" = "  +0.06%  "
    "D19" = "27.250.71"
    "E19" = "  -0.22%  "
    "D20" = "5.031"
    "E20" = "  -1.36%  "
    "D21" = "1.065.86"
    "E21" = "  -43.93%  "
    "D22" = "10.76"
    "E22" = "  +0.07%  "
    "D23" = "6.506"
    "E23" = "  +1.12%  "
    "D24" = "149.64"
    "E24" = "  +1.66%  "
    "D25" = "2.287"
    "E25" = "  -0.72%  "
    "D26" = "18.34"
    "E26" = "  +0.43%  "
    "D27" = "1.746"
    "E27" = "  -0.05%  "
    "D28" = "116.73"
    "E28" = "  +1.04%  "
    "D29" = "4.814"
    "E29" = "  -0.44%  "
    "D30" = "4.783"
    "E30" = "  -2.96%  "
    "D31" = "0.09224"
    "E31" = "  -0.33%  "
    "D32" = "0.8285"
    "E32" = "  +3.64%  "
    "D33" = "0.05057"
    "E33" = "  -0.27%  "
    "D34" = "1.222"
    "E34" = "  -1.21%  "
    "D35" = "3.001"
    "E35" = "  +1.32%  "
    "D36" = "3.338"
    "E36" = "  -3.07%  "
    "D37" = "2.678"
    "E37" = "  +2.85%  "
    "D38" = "0.5743"
    "E38" = "  +0.74%  "
    "D39" = "0.02005"
    "E39" = "  -0.12%  "
    "D40" = "1.075"
    "E40" = "  -0.18%  "
    "D41" = "9.342"
    "E41" = "  +3.71%  "
    "D42" = "6.592"
    "E42" = "  +0.15%  "
    "D43" = "117.02"
    "E43" = "  +1.16%  "
    "E44" = "  +0.07%  "
    "D45" = "0.4947"
    "E45" = "  +0.87%  "
    "E46" = "  -0.06%  "
    "D47" = "10.09"
    "E47" = "  +0.09%  "
    "D48" = "1.636"
    "E48" = "  +0.71%  "
    "D49" = "38.30"
    "E49" = "  +0.39%  "
    "D50" = "0.06174"
    "E50" = "  +3.75%  "
    "D51" = "63.30"
    "E51" = "  -0.85%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}

Write-Output "Updated $($updates.Count) cells"
